$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1794.25
$ws.Range("I40").Value = 1794.25
$ws.Range("K40").Value = 1794.25
$ws.Range("M40").Value = -1619.25

$ws.Range("H62").Value = 10720.8125
$ws.Range("J62").Value = 22501.25
$ws.Range("L62").Value = 22501.25
$ws.Range("N62").Value = -23749.25

$ws.Range("H64").Value = 6500.4
$ws.Range("I64").Value = 3002
$ws.Range("K64").Value = 3002
$ws.Range("M64").Value = -2754

$ws.Range("H65").Value = 10720.8125
$ws.Range("J65").Value = 22501.25
$ws.Range("L65").Value = 112506.25
$ws.Range("N65").Value = -118746.25

$ws.Range("H67").Value = 6500.4
$ws.Range("I67").Value = 3002
$ws.Range("K67").Value = 3002
$ws.Range("M67").Value = -2144

$ws.Range("H76").Value = 6560.6
$ws.Range("I76").Value = 6528.857
$ws.Range("K76").Value = 6528.857
$ws.Range("M76").Value = -6213.857

$ws.Range("H79").Value = 6560.6
$ws.Range("I79").Value = 6528.857
$ws.Range("K79").Value = 6528.857
$ws.Range("M79").Value = -5436.857

$ws.Range("H99").Value = 1536.125
$ws.Range("I99").Value = 1383.1666
$ws.Range("K99").Value = 4149.4998
$ws.Range("M99").Value = -2651.4998

$ws.Range("H118").Value = 4604.5
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 4604.5
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 13813.5
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -17127.5

$ws.Range("H132").Value = 1499.61
$ws.Range("I132").Value = 1505.9584
$ws.Range("K132").Value = 4517.8752
$ws.Range("M132").Value = -1987.8752

$ws.Range("H138").Value = 9777.073
$ws.Range("J138").Value = 10113.926
$ws.Range("L138").Value = 30341.778
$ws.Range("N138").Value = -40621.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1284625.9
$ws.Range("I32").Value = 1523594
$ws.Range("K32").Value = 1523594
$ws.Range("M32").Value = -1523307

$ws.Range("H122").Value = 4697.927
$ws.Range("I122").Value = 3193.6155
$ws.Range("K122").Value = 9580.8465
$ws.Range("M122").Value = -7130.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 29867.75
$ws.Range("I26").Value = 29867.75
$ws.Range("K26").Value = 29867.75
$ws.Range("M26").Value = -29575.75

$ws.Range("H40").Value = 50000
$ws.Range("J40").Value = 50000
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50530

$ws.Range("H105").Value = 3205.7
$ws.Range("I105").Value = 3182.125
$ws.Range("K105").Value = 3182.125
$ws.Range("M105").Value = -1435.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18258.17
$ws.Range("I31").Value = 7638.1714
$ws.Range("J31").Value = 49233.168
$ws.Range("K31").Value = 7638.1714
$ws.Range("L31").Value = 49233.168
$ws.Range("M31").Value = -7343.1714
$ws.Range("N31").Value = -49823.168

$ws.Range("H34").Value = 18258.17
$ws.Range("I34").Value = 7638.1714
$ws.Range("J34").Value = 49233.168
$ws.Range("K34").Value = 7638.1714
$ws.Range("L34").Value = 49233.168
$ws.Range("M34").Value = -7436.1714
$ws.Range("N34").Value = -49637.168

$ws.Range("H107").Value = 481743.3
$ws.Range("I107").Value = 917843.44
$ws.Range("J107").Value = 5997.727
$ws.Range("K107").Value = 917843.44
$ws.Range("L107").Value = 5997.727
$ws.Range("M107").Value = -915923.44
$ws.Range("N107").Value = -9837.726999999999

$ws.Range("H110").Value = 79749.75
$ws.Range("J110").Value = 79749.75
$ws.Range("L110").Value = 79749.75
$ws.Range("N110").Value = -87929.75

$ws.Range("H132").Value = 10162
$ws.Range("I132").Value = 3477.25
$ws.Range("K132").Value = 10431.75
$ws.Range("M132").Value = -7901.75

$ws.Range("H134").Value = 7104.3706
$ws.Range("I134").Value = 1611
$ws.Range("J134").Value = 13971.083
$ws.Range("K134").Value = 4833
$ws.Range("L134").Value = 41913.249
$ws.Range("M134").Value = -2298
$ws.Range("N134").Value = -46983.249

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2508138.2
$ws.Range("I50").Value = 2184.5
$ws.Range("J50").Value = 17543860
$ws.Range("K50").Value = 6553.5
$ws.Range("L50").Value = 52631580
$ws.Range("M50").Value = -6072.5
$ws.Range("N50").Value = -52632542

$ws.Range("H53").Value = 2508138.2
$ws.Range("I53").Value = 2184.5
$ws.Range("J53").Value = 17543860
$ws.Range("K53").Value = 6553.5
$ws.Range("L53").Value = 52631580
$ws.Range("M53").Value = -6072.5
$ws.Range("N53").Value = -52632542

$ws.Range("H55").Value = 2358.1667
$ws.Range("I55").Value = 2429.8
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 7289.400000000001
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = -7112.400000000001
$ws.Range("N55").Value = -6354

$ws.Range("H97").Value = 2500.2173
$ws.Range("J97").Value = 5738.8887
$ws.Range("L97").Value = 17216.6661
$ws.Range("N97").Value = -18208.6661

$ws.Range("H131").Value = 1464.54
$ws.Range("I131").Value = 995.6667
$ws.Range("J131").Value = 1479.0413
$ws.Range("K131").Value = 2987.0001
$ws.Range("L131").Value = 4437.123900000001
$ws.Range("M131").Value = 2052.9999
$ws.Range("N131").Value = -14517.1239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 49995
$ws.Range("J26").Value = 49995
$ws.Range("L26").Value = 49995
$ws.Range("N26").Value = -50555

$ws.Range("H50").Value = 49995
$ws.Range("J50").Value = 49995
$ws.Range("L50").Value = 49995
$ws.Range("N50").Value = -50991

$ws.Range("H70").Value = 8963.75
$ws.Range("I70").Value = 5456.5713
$ws.Range("J70").Value = 10852.23
$ws.Range("K70").Value = 5456.5713
$ws.Range("L70").Value = 10852.23
$ws.Range("M70").Value = -5186.5713
$ws.Range("N70").Value = -11392.23

$ws.Range("H73").Value = 8963.75
$ws.Range("I73").Value = 5456.5713
$ws.Range("J73").Value = 10852.23
$ws.Range("K73").Value = 5456.5713
$ws.Range("L73").Value = 10852.23
$ws.Range("M73").Value = -4520.5713
$ws.Range("N73").Value = -12724.23

$ws.Range("H80").Value = 16713.176
$ws.Range("I80").Value = 9493.444
$ws.Range("K80").Value = 9493.444
$ws.Range("M80").Value = -8495.444

$ws.Range("H83").Value = 16713.176
$ws.Range("I83").Value = 9493.444
$ws.Range("K83").Value = 47467.22
$ws.Range("M83").Value = -42475.22

$ws.Range("H122").Value = 4382.353
$ws.Range("I122").Value = 2950.5
$ws.Range("K122").Value = 8851.5
$ws.Range("M122").Value = -6401.5

$ws.Range("H132").Value = 5700.591
$ws.Range("I132").Value = 5958.3125
$ws.Range("K132").Value = 17874.9375
$ws.Range("M132").Value = -15344.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3200.5
$ws.Range("I16").Value = 3050.6155
$ws.Range("J16").Value = 3850
$ws.Range("K16").Value = 3050.6155
$ws.Range("L16").Value = 3850
$ws.Range("M16").Value = -2880.6155
$ws.Range("N16").Value = -4190

$ws.Range("H21").Value = 1500
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H46").Value = 2937.125
$ws.Range("I46").Value = 3499.3333
$ws.Range("K46").Value = 3499.3333
$ws.Range("M46").Value = -3311.3333

$ws.Range("H55").Value = 1565.6522
$ws.Range("I55").Value = 935.1429
$ws.Range("J55").Value = 1841.5
$ws.Range("K55").Value = 935.1429
$ws.Range("L55").Value = 1841.5
$ws.Range("M55").Value = -762.1429
$ws.Range("N55").Value = -2187.5

$ws.Range("H122").Value = 9417.174
$ws.Range("I122").Value = 6773.067
$ws.Range("K122").Value = 20319.201
$ws.Range("M122").Value = -17869.201

$ws.Range("H136").Value = 15259.86
$ws.Range("I136").Value = 13045.192
$ws.Range("K136").Value = 39135.576
$ws.Range("M136").Value = -36585.576

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 879.6667
$ws.Range("I100").Value = 817.8889
$ws.Range("K100").Value = 1635.7778
$ws.Range("M100").Value = -1094.7778

$ws.Range("H136").Value = 8743.206
$ws.Range("I136").Value = 2220.5
$ws.Range("K136").Value = 6661.5
$ws.Range("M136").Value = -4111.5

